# Fruta / hortaliza, semanal
# Insert two new weekly data rows at the top of the Cilantro data block
# (row 1232), pushing the existing rows 1232:1335 down to 1234:1337.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 1232 - this shifts rows 1232:1335 down to 1234:1337
# and Excel copies the formatting (incl. the date number format) from the
# row above, same as a normal "Insert Copied Cells"/"Insert Rows" in the UI.
$ws.Rows("1232:1233").Insert()

# --- New row 1232 ---
$ws.Range("A1232").Value = 6
$ws.Range("B1232").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1232").Value = "Metropolitana"
$ws.Range("D1232").Value = 45106
$ws.Range("E1232").Value = 13
$ws.Range("F1232").Value = 100112040
$ws.Range("G1232").Value = "Cilantro"
$ws.Range("H1232").Value = "Sin especificar"
$ws.Range("I1232").Value = "Primera"
$ws.Range("J1232").Value = 710
$ws.Range("K1232").Value = 6000
$ws.Range("L1232").Value = 7000
$ws.Range("M1232").Value = 6465
$ws.Range("N1232").Value = "`$/caja 36 atados"
$ws.Range("O1232").Value = "Región Metropolitana"
$ws.Range("P1232").Value = 180
$ws.Range("Q1232").Value = 36
$ws.Range("R1232").Value = "Hortaliza"

# --- New row 1233 ---
$ws.Range("A1233").Value = 6
$ws.Range("B1233").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1233").Value = "Metropolitana"
$ws.Range("D1233").Value = 45106
$ws.Range("E1233").Value = 13
$ws.Range("F1233").Value = 100112040
$ws.Range("G1233").Value = "Cilantro"
$ws.Range("H1233").Value = "Sin especificar"
$ws.Range("I1233").Value = "Primera"
$ws.Range("J1233").Value = 440
$ws.Range("K1233").Value = 14000
$ws.Range("L1233").Value = 15000
$ws.Range("M1233").Value = 14409
$ws.Range("N1233").Value = "`$/docena de atados"
$ws.Range("O1233").Value = "Región Metropolitana"
$ws.Range("P1233").Value = 4803
$ws.Range("Q1233").Value = 3
$ws.Range("R1233").Value = "Hortaliza"
